$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 152.72728
$ws.Range("I2").Value = 76.77778000000001
$ws.Range("J2").Value = 494.5
$ws.Range("K2").Value = 76.77778000000001
$ws.Range("L2").Value = 494.5
$ws.Range("M2").Value = 36.22221999999999
$ws.Range("N2").Value = -720.5
$ws.Range("H6").Value = 38.210526
$ws.Range("I6").Value = 38.210526
$ws.Range("K6").Value = 114.631578
$ws.Range("M6").Value = -2.631578000000005
$ws.Range("H40").Value = 781.25
$ws.Range("I40").Value = 546.875
$ws.Range("K40").Value = 546.875
$ws.Range("M40").Value = -371.875
$ws.Range("H74").Value = 2966.2778
$ws.Range("I74").Value = 2399.5625
$ws.Range("K74").Value = 2399.5625
$ws.Range("M74").Value = -1463.5625
$ws.Range("H77").Value = 2966.2778
$ws.Range("I77").Value = 2399.5625
$ws.Range("K77").Value = 11997.8125
$ws.Range("M77").Value = -7317.8125
$ws.Range("H116").Value = 5998.4707
$ws.Range("J116").Value = 5758.4
$ws.Range("L116").Value = 5758.4
$ws.Range("N116").Value = -12642.4
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H137").Value = 34286.43
$ws.Range("I137").Value = 11999
$ws.Range("K137").Value = 35997
$ws.Range("M137").Value = -33447
$ws.Range("H138").Value = 4278.404
$ws.Range("J138").Value = 4487.125
$ws.Range("L138").Value = 13461.375
$ws.Range("N138").Value = -23741.375
$ws.Range("H141").Value = 3748.05
$ws.Range("I141").Value = 3637.4119
$ws.Range("J141").Value = 4375
$ws.Range("K141").Value = 10912.2357
$ws.Range("L141").Value = 13125
$ws.Range("M141").Value = -5732.235700000001
$ws.Range("N141").Value = -23485

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2694.7334
$ws.Range("I102").Value = 2530.4285
$ws.Range("J102").Value = 4995
$ws.Range("K102").Value = 2530.4285
$ws.Range("L102").Value = 4995
$ws.Range("M102").Value = -908.4285
$ws.Range("N102").Value = -8239
$ws.Range("H122").Value = 2303.125
$ws.Range("I122").Value = 2189.2856
$ws.Range("K122").Value = 6567.8568
$ws.Range("M122").Value = -4117.8568
$ws.Range("H125").Value = 49500
$ws.Range("J125").Value = 49500
$ws.Range("L125").Value = 49500
$ws.Range("N125").Value = -59340
$ws.Range("H132").Value = 16548261
$ws.Range("I132").Value = 5602
$ws.Range("K132").Value = 16806
$ws.Range("M132").Value = -14276

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 995
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H54").Value = 6488.25
$ws.Range("J54").Value = 10000
$ws.Range("L54").Value = 10000
$ws.Range("N54").Value = -10968

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10789.182
$ws.Range("I31").Value = 1442.9375
$ws.Range("J31").Value = 35712.5
$ws.Range("K31").Value = 1442.9375
$ws.Range("L31").Value = 35712.5
$ws.Range("M31").Value = -1147.9375
$ws.Range("N31").Value = -36302.5
$ws.Range("H34").Value = 10789.182
$ws.Range("I34").Value = 1442.9375
$ws.Range("J34").Value = 35712.5
$ws.Range("K34").Value = 1442.9375
$ws.Range("L34").Value = 35712.5
$ws.Range("M34").Value = -1240.9375
$ws.Range("N34").Value = -36116.5
$ws.Range("H58").Value = 24011.5
$ws.Range("I58").Value = 13016.625
$ws.Range("K58").Value = 13016.625
$ws.Range("M58").Value = -12813.625
$ws.Range("H93").Value = 9333
$ws.Range("I93").Value = 9333
$ws.Range("K93").Value = 9333
$ws.Range("M93").Value = -7461
$ws.Range("H132").Value = 42894492
$ws.Range("I132").Value = 2573.1904
$ws.Range("J132").Value = 343137920
$ws.Range("K132").Value = 7719.5712
$ws.Range("L132").Value = 1029413760
$ws.Range("M132").Value = -5189.5712
$ws.Range("N132").Value = -1029418820
$ws.Range("H136").Value = 24011.5
$ws.Range("I136").Value = 13016.625
$ws.Range("K136").Value = 39049.875
$ws.Range("M136").Value = -36499.875

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 5172.769
$ws.Range("I107").Value = 841
$ws.Range("J107").Value = 8885.714
$ws.Range("K107").Value = 2523
$ws.Range("L107").Value = 26657.142
$ws.Range("M107").Value = -603
$ws.Range("N107").Value = -30497.142

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 2011980
$ws.Range("I20").Value = 5007450
$ws.Range("J20").Value = 15000
$ws.Range("K20").Value = 5007450
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = -5007205
$ws.Range("N20").Value = -15490
$ws.Range("H24").Value = 2233889
$ws.Range("J24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("N24").Value = -15346
$ws.Range("H113").Value = 3862
$ws.Range("J113").Value = 4462.375
$ws.Range("L113").Value = 4462.375
$ws.Range("N113").Value = -8802.375
$ws.Range("H126").Value = 8223.182000000001
$ws.Range("I126").Value = 10158.571
$ws.Range("J126").Value = 4836.25
$ws.Range("K126").Value = 30475.713
$ws.Range("L126").Value = 14508.75
$ws.Range("M126").Value = -28005.713
$ws.Range("N126").Value = -19448.75
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H35").Value = 2516.2
$ws.Range("I35").Value = 1677
$ws.Range("K35").Value = 1677
$ws.Range("M35").Value = -1341
$ws.Range("H40").Value = 4663.4443
$ws.Range("J40").Value = 4499.5
$ws.Range("L40").Value = 4499.5
$ws.Range("N40").Value = -4771.5
$ws.Range("H46").Value = 2605.5518
$ws.Range("I46").Value = 1160
$ws.Range("J46").Value = 3366.3684
$ws.Range("K46").Value = 1160
$ws.Range("L46").Value = 3366.3684
$ws.Range("M46").Value = -972
$ws.Range("N46").Value = -3742.3684
$ws.Range("H119").Value = 75000
$ws.Range("J119").Value = 75000
$ws.Range("L119").Value = 75000
$ws.Range("N119").Value = -84676
$ws.Range("H122").Value = 10203.75
$ws.Range("I122").Value = 11407.5
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 34222.5
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -31772.5
$ws.Range("N122").Value = -31900
$ws.Range("H132").Value = 4988970.5
$ws.Range("I132").Value = 3692.2222
$ws.Range("K132").Value = 11076.6666
$ws.Range("M132").Value = -8546.6666
$ws.Range("H136").Value = 155521.53
$ws.Range("I136").Value = 22181.5
$ws.Range("K136").Value = 66544.5
$ws.Range("M136").Value = -63994.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 2680000
$ws.Range("J18").Value = 4015000
$ws.Range("L18").Value = 4015000
$ws.Range("N18").Value = -4015346
$ws.Range("H21").Value = 38194.09
$ws.Range("I21").Value = 30015
$ws.Range("K21").Value = 30015
$ws.Range("M21").Value = -29780
$ws.Range("H22").Value = 1995
$ws.Range("I22").Value = 1995
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1995
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1702
$ws.Range("N22").ClearContents()
$ws.Range("H35").Value = 38194.09
$ws.Range("I35").Value = 30015
$ws.Range("K35").Value = 30015
$ws.Range("M35").Value = -29725
$ws.Range("H107").Value = 1014.069
$ws.Range("I107").Value = 1066.08
$ws.Range("K107").Value = 3198.24
$ws.Range("M107").Value = -1278.24
$ws.Range("H126").Value = 4164.6665
$ws.Range("I126").Value = 4310.25
$ws.Range("K126").Value = 12930.75
$ws.Range("M126").Value = -10460.75
$ws.Range("H132").Value = 1210161.8
$ws.Range("I132").Value = 2484.4285
$ws.Range("J132").Value = 5437032.5
$ws.Range("K132").Value = 7453.2855
$ws.Range("L132").Value = 16311097.5
$ws.Range("M132").Value = -4923.2855
$ws.Range("N132").Value = -16316157.5
